$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Author FullName" header text to "Author Full Name"
$ws.Range("A1").Value = "Author Full Name"

# Replace the author value "haruki murakami" with "theodore dreiser"
$ws.Range("A3").Value = "theodore dreiser"

# Set the page orientation to portrait (adds <pageSetup orientation="portrait".../>)
$ws.PageSetup.Orientation = 1

# Move the active selection to G11 (updates <selection activeCell="G11" sqref="G11"/>)
$ws.Range("G11").Select()
